$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.62
$ws.Range("H3").Value = 3.6
$ws.Range("I3").Value = 2.42
$ws.Range("J3").Value = 1.05
$ws.Range("K3").Value = 8.25
$ws.Range("L3").Value = 1.27
$ws.Range("M3").Value = 3.5
$ws.Range("N3").Value = 1.82
$ws.Range("O3").Value = 1.93
$ws.Range("P3").Value = 1.38
$ws.Range("Q3").Value = 2.92
$ws.Range("R3").Value = 1.7
$ws.Range("S3").Value = 2.05
$ws.Range("T3").Value = 9
$ws.Range("U3").Value = 14.5
$ws.Range("V3").Value = 10.5
$ws.Range("W3").Value = 30
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 32
$ws.Range("Z3").Value = 8.25
$ws.Range("AA3").Value = 7.2
$ws.Range("AB3").Value = 14.5
$ws.Range("AC3").Value = 65
$ws.Range("AD3").Value = 500
$ws.Range("AE3").Value = 8.75
$ws.Range("AF3").Value = 13
$ws.Range("AG3").Value = 10
$ws.Range("AH3").Value = 27
$ws.Range("AI3").Value = 21
$ws.Range("AJ3").Value = 30

# Row 9
$ws.Range("G9").Value = 1.95
$ws.Range("H9").Value = 4.1
$ws.Range("I9").Value = 3
$ws.Range("N9").Value = 1.35
$ws.Range("O9").Value = 2.72
$ws.Range("R9").Value = 1.36
$ws.Range("S9").Value = 2.65
$ws.Range("T9").Value = 14.5
$ws.Range("U9").Value = 14
$ws.Range("V9").Value = 9.25
$ws.Range("W9").Value = 20
$ws.Range("X9").Value = 13.5
$ws.Range("Y9").Value = 16.5
$ws.Range("Z9").Value = 23
$ws.Range("AA9").Value = 9.5
$ws.Range("AB9").Value = 11.5
$ws.Range("AC9").Value = 30
$ws.Range("AD9").Value = 120
$ws.Range("AE9").Value = 18
$ws.Range("AF9").Value = 22
$ws.Range("AG9").Value = 11.75
$ws.Range("AH9").Value = 40
$ws.Range("AI9").Value = 21
$ws.Range("AJ9").Value = 21

# Row 10
$ws.Range("G10").Value = 2.02
$ws.Range("H10").Value = 3.85
$ws.Range("I10").Value = 2.95
$ws.Range("N10").Value = 1.5
$ws.Range("O10").Value = 2.25
$ws.Range("R10").Value = 1.47
$ws.Range("S10").Value = 2.32
$ws.Range("T10").Value = 11
$ws.Range("U10").Value = 12
$ws.Range("V10").Value = 8.75
$ws.Range("W10").Value = 19.5
$ws.Range("X10").Value = 14.5
$ws.Range("Y10").Value = 20
$ws.Range("Z10").Value = 16.5
$ws.Range("AA10").Value = 7.9
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 40
$ws.Range("AD10").Value = 250
$ws.Range("AE10").Value = 13.5
$ws.Range("AF10").Value = 18.5
$ws.Range("AG10").Value = 11
$ws.Range("AH10").Value = 37
$ws.Range("AI10").Value = 22
$ws.Range("AJ10").Value = 24

# Row 13
$ws.Range("G13").Value = 3.2
$ws.Range("H13").Value = 3.05
$ws.Range("I13").Value = 2.22
$ws.Range("L13").Value = 1.45
$ws.Range("M13").Value = 2.37
$ws.Range("Q13").Value = 2.27
$ws.Range("R13").Value = 2.02
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 7.6
$ws.Range("U13").Value = 15
$ws.Range("V13").Value = 12
$ws.Range("W13").Value = 45
$ws.Range("X13").Value = 35
$ws.Range("Z13").Value = 6.9
$ws.Range("AA13").Value = 6.1
$ws.Range("AB13").Value = 18
$ws.Range("AC13").Value = 110
$ws.Range("AE13").Value = 5.9
$ws.Range("AF13").Value = 9.25
$ws.Range("AG13").Value = 9.5
$ws.Range("AH13").Value = 21
$ws.Range("AI13").Value = 22

# Row 14
$ws.Range("G14").Value = 3.35
$ws.Range("H14").Value = 3.05
$ws.Range("I14").Value = 2.15
$ws.Range("L14").Value = 1.5
$ws.Range("M14").Value = 2.27
$ws.Range("N14").Value = 2.42
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 1.52
$ws.Range("Q14").Value = 2.22
$ws.Range("R14").Value = 2.1
$ws.Range("S14").Value = 1.57
$ws.Range("T14").Value = 7.4
$ws.Range("U14").Value = 15.5
$ws.Range("V14").Value = 13
$ws.Range("W14").Value = 45
$ws.Range("X14").Value = 40
$ws.Range("Y14").Value = 60
$ws.Range("Z14").Value = 6.5
$ws.Range("AA14").Value = 6.1
$ws.Range("AB14").Value = 20
$ws.Range("AC14").Value = 150
$ws.Range("AE14").Value = 5.6
$ws.Range("AF14").Value = 9
$ws.Range("AG14").Value = 9.75
$ws.Range("AH14").Value = 20
$ws.Range("AI14").Value = 22
$ws.Range("AJ14").Value = 45

# Row 16
$ws.Range("G16").Value = 2.67
$ws.Range("H16").Value = 3.45
$ws.Range("I16").Value = 2.3
